$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto table: "Price" (column D) and "Volume(1h)" (column E)
# for each coin row, as pulled from the latest coinranking.com snapshot.
#
# Some new values (e.g. "1.00", "0.800", "7.60", "0.0000239") look like plain
# numbers and would otherwise be auto-converted/normalized by Excel, dropping
# trailing zeros or switching to scientific notation. To keep them identical
# to the source text (and keep the "." used as a thousands separator in big
# prices, e.g. "93.440.06"), each cell is temporarily forced to text format,
# written, and then restored to its original style so no formatting metadata
# changes are introduced.
$updates = [ordered]@{
    "D2" = "93.440.06"
    "E2" = "  -0.28%  "
    "D3" = "3.048.14"
    "E3" = "  -2.41%  "
    "E4" = "  +0.09%  "
    "D5" = "233.42"
    "E5" = "  -4.15%  "
    "D6" = "605.92"
    "E6" = "  -1.82%  "
    "E7" = "  +0.05%  "
    "D8" = "0.374"
    "E8" = "  -9.63%  "
    "D9" = "1.00"
    "E9" = "  +0.12%  "
    "D10" = "0.800"
    "E10" = "  +8.52%  "
    "D11" = "3.048.23"
    "E11" = "  -2.38%  "
    "D12" = "0.195"
    "E12" = "  -4.24%  "
    "D13" = "93.373.67"
    "E13" = "  +0.52%  "
    "D14" = "0.0000239"
    "E14" = "  -7.26%  "
    "D15" = "33.44"
    "E15" = "  -3.82%  "
    "D16" = "5.26"
    "E16" = "  -4.28%  "
    "D17" = "3.624.14"
    "E17" = "  -2.36%  "
    "D18" = "3.036.61"
    "E18" = "  -2.70%  "
    "D19" = "3.53"
    "E19" = "  -7.24%  "
    "D20" = "14.33"
    "E20" = "  -3.27%  "
    "D21" = "5.67"
    "E21" = "  -2.67%  "
    "D22" = "436.18"
    "E22" = "  -3.51%  "
    "D23" = "8.75"
    "E23" = "  -7.44%  "
    "D24" = "0.0000189"
    "E24" = "  -9.93%  "
    "D25" = "8.33"
    "E25" = "  +5.19%  "
    "D26" = "5.47"
    "E26" = "  -6.41%  "
    "D27" = "84.04"
    "E27" = "  -4.04%  "
    "D28" = "11.72"
    "E28" = "  -1.34%  "
    "D29" = "3.223.74"
    "E29" = "  -2.08%  "
    "E30" = "  +0.06%  "
    "D31" = "0.248"
    "E31" = "  +9.83%  "
    "E32" = "  +4.13%  "
    "D33" = "0.121"
    "E33" = "  -10.49%  "
    "E34" = "  -0.07%  "
    "D35" = "9.02"
    "E35" = "  -2.47%  "
    "D36" = "7.60"
    "E36" = "  -6.52%  "
    "D37" = "0.155"
    "E37" = "  -6.07%  "
    "D38" = "25.26"
    "E38" = "  -3.85%  "
    "D39" = "1.87"
    "E39" = "  -2.53%  "
    "D40" = "23.91"
    "E40" = "  +3.44%  "
    "D41" = "0.438"
    "E41" = "  +0.15%  "
    "D42" = "3.74"
    "E42" = "  -6.11%  "
    "D43" = "459.73"
    "E43" = "  -4.87%  "
    "D44" = "1.24"
    "E44" = "  -5.31%  "
    "D46" = "3.11"
    "E46" = "  -11.86%  "
    "D47" = "160.04"
    "E47" = "  -1.45%  "
    "D48" = "1.82"
    "E48" = "  -6.42%  "
    "D49" = "0.663"
    "E49" = "  -4.62%  "
    "E50" = "  -0.83%  "
    "E51" = "  +0.04%  "
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = $origStyle
}
